{"js": "// Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific bullet /\n// impact paragraphs, splitting the existing single run into multiple runs\n// exactly as produced by the canonical diff.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Paragraphs to update, identified by their exact current text, together\n// with the list of substrings (in left-to-right order, each occurring\n// exactly once) that must become bold + colored.\nconst targets = [\n  {\n    match:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    terms: [\"23%\", \"64%\"],\n  },\n  {\n    match:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    terms: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    match:\n      \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    terms: [\"1,200\"],\n  },\n  {\n    match:\n      \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    terms: [\"$400M\", \"$1B\"],\n  },\n  {\n    match:\n      \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    terms: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    match:\n      \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    terms: [\"87%\", \"71%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (const target of targets) {\n  const para = paragraphs.items.find((p) => p.text === target.match);\n  if (!para) {\n    continue;\n  }\n  const paraRange = para.getRange();\n\n  for (const term of target.terms) {\n    const searchResults = paraRange.search(term, { matchCase: true });\n    searchResults.load(\"items\");\n    await context.sync();\n\n    if (searchResults.items.length === 0) {\n      continue;\n    }\n    const hit = searchResults.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific bullet /\n# impact paragraphs, splitting the existing single run into multiple runs\n# exactly as produced by the canonical diff.\n\n$d = $word.ActiveDocument\n\n# RGB(0x2C,0x3E,0x50) expressed as the BGR integer Word's\n# Range.Font.Color property expects.\n$highlightColor = 0x2C + (0x3E * 256) + (0x50 * 65536)\n\n$bullet = [char]0x2022\n$plusMinus = [char]0x00B1\n\nfunction Get-ParagraphByText($doc, $targetText) {\n    foreach ($para in $doc.Paragraphs) {\n        $raw = $para.Range.Text\n        $trimmed = $raw.TrimEnd([char]13, [char]7)\n        if ($trimmed -eq $targetText) {\n            return $para\n        }\n    }\n    return $null\n}\n\nfunction Highlight-Term($scopeRange, $text) {\n    $searchRange = $scopeRange.Duplicate()\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $text\n    $find.MatchCase = $true\n    $found = $find.Execute()\n    if ($found) {\n        $searchRange.Bold = 1\n        $searchRange.Font.Color = $highlightColor\n    }\n    return $found\n}\n\n$targets = @(\n    @{\n        Match = $bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Terms = @(\"23%\", \"64%\")\n    },\n    @{\n        Match = $bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $plusMinus + \"4.2% to \" + $plusMinus + \"2.1%\"\n        Terms = @(\"87%\", \"71%\", ($plusMinus + \"4.2%\"), ($plusMinus + \"2.1%\"))\n    },\n    @{\n        Match = $bullet + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Terms = @(\"1,200\")\n    },\n    @{\n        Match = $bullet + \" Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Terms = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Match = $bullet + \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Terms = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Match = $bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Terms = @(\"87%\", \"71%\")\n    }\n)\n\nforeach ($target in $targets) {\n    $para = Get-ParagraphByText $d $target.Match\n    if ($para -eq $null) {\n        continue\n    }\n    $paraRange = $para.Range.Duplicate()\n    foreach ($term in $target.Terms) {\n        Highlight-Term $paraRange $term | Out-Null\n    }\n}\n"}
